# Add a 2022-Q3 sheet (with fund-holding detail) and a summary row in "总计".
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for 2022-Q3, shift the rest down, and
#    renumber the running index in column A (0,1,2,3,4).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# New row 2 gets column-A style matching the other index cells (copy format
# only, so the freshly-set numeric value below survives).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 7
$summary.Cells.Item(2,4).Value = 0.36

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: duplicate the "2022-Q2" sheet (so headers,
#    column widths and formatting all match), place it right before
#    "2022-Q2", rename it, and overwrite its data with the Q3 fund list.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2) | Out-Null
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Make sure there are 8 rows (1 header + 7 funds) with matching formatting:
# clone row 3's look (index + plain-text cells) down through row 8.
$q3.Range("A3:H3").Copy()
$q3.Range("A4:H8").PasteSpecial(-4122)

# Fund code / size / position figures are text in the source data (keeps
# leading/trailing zeros like "013331" or "8.10") - force that explicitly.
$q3.Range("B2:B8").NumberFormat = "@"
$q3.Range("D2:G8").NumberFormat = "@"

$q3Data = @(
    @(0, "161039", "富国中证1000指数增强（LOF）A", "25.41", "84.72", "0.63", "0.1601", 7),
    @(1, "013331", "富国中证1000指数增强（LOF）C", "8.53",  "84.72", "0.63", "0.0537", 7),
    @(2, "015784", "中信建投中证1000指数增强A",     "8.10",  "92.20", "0.64", "0.0518", 9),
    @(3, "006165", "建信中证1000指数增强A",         "3.87",  "84.02", "1.18", "0.0457", 10),
    @(4, "006166", "建信中证1000指数增强C",         "1.89",  "84.02", "1.18", "0.0223", 10),
    @(5, "015785", "中信建投中证1000指数增强C",     "3.32",  "92.20", "0.64", "0.0212", 9),
    @(6, "013442", "建信中证1000指数增强E",         "0.18",  "84.02", "1.18", "0.0021", 10)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r,1).Value = $row[0]
    $q3.Cells.Item($r,2).Value = $row[1]
    $q3.Cells.Item($r,3).Value = $row[2]
    $q3.Cells.Item($r,4).Value = $row[3]
    $q3.Cells.Item($r,5).Value = $row[4]
    $q3.Cells.Item($r,6).Value = $row[5]
    $q3.Cells.Item($r,7).Value = $row[6]
    $q3.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}
